$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.024338149013958
$ws.Cells.Item(2, 4).Value = 1.034191466648898
$ws.Cells.Item(2, 5).Value = 1.024815386143875
$ws.Cells.Item(2, 6).Value = 1.022818776234388
$ws.Cells.Item(2, 9).Value = 1.033440287926734
$ws.Cells.Item(2, 10).Value = 1.029513323041823
$ws.Cells.Item(2, 11).Value = 1.036991649358475
$ws.Cells.Item(2, 12).Value = 1.027642781320001
$ws.Cells.Item(2, 13).Value = 1.025652036322532
$ws.Cells.Item(2, 14).Value = 1.013800290477456
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025546527427986
$ws.Cells.Item(3, 4).Value = 1.034785196301982
$ws.Cells.Item(3, 5).Value = 1.025848877626426
$ws.Cells.Item(3, 6).Value = 1.024667743269559
$ws.Cells.Item(3, 9).Value = 1.033681192818823
$ws.Cells.Item(3, 10).Value = 1.03035922663561
$ws.Cells.Item(3, 11).Value = 1.037395090002931
$ws.Cells.Item(3, 12).Value = 1.028482828766609
$ws.Cells.Item(3, 13).Value = 1.027304908121812
$ws.Cells.Item(3, 14).Value = 1.014087431925886
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.026327410116028
$ws.Cells.Item(4, 4).Value = 1.035169016740766
$ws.Cells.Item(4, 5).Value = 1.026517062960678
$ws.Cells.Item(4, 6).Value = 1.02586283580272
$ws.Cells.Item(4, 9).Value = 1.033835530653418
$ws.Cells.Item(4, 10).Value = 1.030905096233534
$ws.Cells.Item(4, 11).Value = 1.037655084728173
$ws.Cells.Item(4, 12).Value = 1.029025247975916
$ws.Cells.Item(4, 13).Value = 1.02837270958878
$ws.Cells.Item(4, 14).Value = 1.014272551847992
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.026655453556383
$ws.Cells.Item(5, 4).Value = 1.035330287874758
$ws.Cells.Item(5, 5).Value = 1.026797838620939
$ws.Cells.Item(5, 6).Value = 1.026364950408374
$ws.Cells.Item(5, 9).Value = 1.033900045346491
$ws.Cells.Item(5, 10).Value = 1.031134227043078
$ws.Cells.Item(5, 11).Value = 1.03776413335433
$ws.Cells.Item(5, 12).Value = 1.029253009026966
$ws.Cells.Item(5, 13).Value = 1.028821212901431
$ws.Cells.Item(5, 14).Value = 1.014350214441138
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.02671051956514
$ws.Cells.Item(6, 4).Value = 1.035357360900291
$ws.Cells.Item(6, 5).Value = 1.026844974581717
$ws.Cells.Item(6, 6).Value = 1.026449240280346
$ws.Cells.Item(6, 9).Value = 1.03391085603144
$ws.Cells.Item(6, 10).Value = 1.031172678517985
$ws.Cells.Item(6, 11).Value = 1.03778242825929
$ws.Cells.Item(6, 12).Value = 1.029291235240353
$ws.Cells.Item(6, 13).Value = 1.028896495408565
$ws.Cells.Item(6, 14).Value = 1.01436324486957
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.026331794384854
$ws.Cells.Item(7, 4).Value = 1.035171171995648
$ws.Cells.Item(7, 5).Value = 1.026520815203918
$ws.Cells.Item(7, 6).Value = 1.02586954625568
$ws.Cells.Item(7, 9).Value = 1.033836394150777
$ws.Cells.Item(7, 10).Value = 1.030908159272928
$ws.Cells.Item(7, 11).Value = 1.037656542835748
$ws.Cells.Item(7, 12).Value = 1.029028292395539
$ws.Cells.Item(7, 13).Value = 1.02837870407225
$ws.Cells.Item(7, 14).Value = 1.014273590214082
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.024746739730342
$ws.Cells.Item(8, 4).Value = 1.03439219525742
$ws.Cells.Item(8, 5).Value = 1.025164775621776
$ws.Cells.Item(8, 6).Value = 1.023443922353071
$ws.Cells.Item(8, 9).Value = 1.033522022893316
$ws.Cells.Item(8, 10).Value = 1.029799509817373
$ws.Cells.Item(8, 11).Value = 1.03712821328442
$ws.Cells.Item(8, 12).Value = 1.027926918590949
$ws.Cells.Item(8, 13).Value = 1.026210994056229
$ws.Cells.Item(8, 14).Value = 1.013897472588348
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021945673506869
$ws.Cells.Item(9, 4).Value = 1.033016772599212
$ws.Cells.Item(9, 5).Value = 1.02277089707642
$ws.Cells.Item(9, 6).Value = 1.019159070455271
$ws.Cells.Item(9, 9).Value = 1.032956208059162
$ws.Cells.Item(9, 10).Value = 1.027834401391574
$ws.Cells.Item(9, 11).Value = 1.036189112857879
$ws.Cells.Item(9, 12).Value = 1.025977242538397
$ws.Cells.Item(9, 13).Value = 1.022377582296327
$ws.Cells.Item(9, 14).Value = 1.013229457092732
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020072645843813
$ws.Cells.Item(10, 4).Value = 1.032097973420878
$ws.Cells.Item(10, 5).Value = 1.02117185951231
$ws.Cells.Item(10, 6).Value = 1.01629462639839
$ws.Cells.Item(10, 9).Value = 1.032570986672975
$ws.Cells.Item(10, 10).Value = 1.026516379885352
$ws.Cells.Item(10, 11).Value = 1.035557571653186
$ws.Cells.Item(10, 12).Value = 1.024671287534228
$ws.Cells.Item(10, 13).Value = 1.019812139631722
$ws.Cells.Item(10, 14).Value = 1.012780525353227
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.019260202449427
$ws.Cells.Item(11, 4).Value = 1.03169968517154
$ws.Cells.Item(11, 5).Value = 1.020478676931121
$ws.Cells.Item(11, 6).Value = 1.015052254170765
$ws.Cells.Item(11, 9).Value = 1.03240227138393
$ws.Cells.Item(11, 10).Value = 1.025943731429865
$ws.Cells.Item(11, 11).Value = 1.035282804263353
$ws.Cells.Item(11, 12).Value = 1.024104290730036
$ws.Cells.Item(11, 13).Value = 1.018698792205592
$ws.Cells.Item(11, 14).Value = 1.012585268332793
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.018958207550329
$ws.Cells.Item(12, 4).Value = 1.031551676723379
$ws.Cells.Item(12, 5).Value = 1.020221075839792
$ws.Cells.Item(12, 6).Value = 1.014590459421493
$ws.Cells.Item(12, 9).Value = 1.032339314850109
$ws.Cells.Item(12, 10).Value = 1.025730729385057
$ws.Cells.Item(12, 11).Value = 1.035180546736526
$ws.Cells.Item(12, 12).Value = 1.023893452507349
$ws.Cells.Item(12, 13).Value = 1.01828485777102
$ws.Cells.Item(12, 14).Value = 1.01251260977158
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.019022996399512
$ws.Cells.Item(13, 4).Value = 1.03158342804605
$ws.Cells.Item(13, 5).Value = 1.020276337738921
$ws.Cells.Item(13, 6).Value = 1.014689530786259
$ws.Cells.Item(13, 9).Value = 1.032352832295295
$ws.Cells.Item(13, 10).Value = 1.02577643245641
$ws.Cells.Item(13, 11).Value = 1.035202490231576
$ws.Cells.Item(13, 12).Value = 1.023938688503701
$ws.Cells.Item(13, 13).Value = 1.018373665862602
$ws.Cells.Item(13, 14).Value = 1.01252820124577
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.019235243920076
$ws.Cells.Item(14, 4).Value = 1.031687452103037
$ws.Cells.Item(14, 5).Value = 1.020457386059353
$ws.Cells.Item(14, 6).Value = 1.015014088733885
$ws.Cells.Item(14, 9).Value = 1.032397073261042
$ws.Cells.Item(14, 10).Value = 1.025926130653664
$ws.Cells.Item(14, 11).Value = 1.035274355639482
$ws.Cells.Item(14, 12).Value = 1.024086867496813
$ws.Cells.Item(14, 13).Value = 1.018664584238967
$ws.Cells.Item(14, 14).Value = 1.012579265041234
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.019365987672622
$ws.Cells.Item(15, 4).Value = 1.031751535950062
$ws.Cells.Item(15, 5).Value = 1.020568919598278
$ws.Cells.Item(15, 6).Value = 1.015214016440656
$ws.Cells.Item(15, 9).Value = 1.032424293367868
$ws.Cells.Item(15, 10).Value = 1.026018325455284
$ws.Cells.Item(15, 11).Value = 1.035318608211625
$ws.Cells.Item(15, 12).Value = 1.024178134832627
$ws.Cells.Item(15, 13).Value = 1.018843776898042
$ws.Cells.Item(15, 14).Value = 1.012610709681382
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020126534867729
$ws.Cells.Item(16, 4).Value = 1.032124397182593
$ws.Cells.Item(16, 5).Value = 1.021217846822788
$ws.Cells.Item(16, 6).Value = 1.01637703414199
$ws.Cells.Item(16, 9).Value = 1.032582143391213
$ws.Cells.Item(16, 10).Value = 1.026554343539633
$ws.Cells.Item(16, 11).Value = 1.035575779505039
$ws.Cells.Item(16, 12).Value = 1.024708885135258
$ws.Cells.Item(16, 13).Value = 1.01988597518793
$ws.Cells.Item(16, 14).Value = 1.012793465556978
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.020603224716647
$ws.Cells.Item(17, 4).Value = 1.032358164792607
$ws.Cells.Item(17, 5).Value = 1.021624687941203
$ws.Cells.Item(17, 6).Value = 1.017106005603043
$ws.Cells.Item(17, 9).Value = 1.03268064602552
$ws.Cells.Item(17, 10).Value = 1.026890052281045
$ws.Cells.Item(17, 11).Value = 1.035736746300398
$ws.Cells.Item(17, 12).Value = 1.02504140409843
$ws.Cells.Item(17, 13).Value = 1.020539041857351
$ws.Cells.Item(17, 14).Value = 1.012907870621186
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.020881134158711
$ws.Cells.Item(18, 4).Value = 1.032494474828481
$ws.Cells.Item(18, 5).Value = 1.021861915631773
$ws.Cells.Item(18, 6).Value = 1.017531005607962
$ws.Cells.Item(18, 9).Value = 1.032737916539735
$ws.Cells.Item(18, 10).Value = 1.027085678942848
$ws.Cells.Item(18, 11).Value = 1.035830509544976
$ws.Cells.Item(18, 12).Value = 1.025235211365198
$ws.Cells.Item(18, 13).Value = 1.02091972485986
$ws.Cells.Item(18, 14).Value = 1.012974517626597
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.020975871245079
$ws.Cells.Item(19, 4).Value = 1.032540945793212
$ws.Cells.Item(19, 5).Value = 1.021942791389603
$ws.Cells.Item(19, 6).Value = 1.017675886698357
$ws.Cells.Item(19, 9).Value = 1.032757413054247
$ws.Cells.Item(19, 10).Value = 1.027152351053646
$ws.Cells.Item(19, 11).Value = 1.03586245903698
$ws.Cells.Item(19, 12).Value = 1.025301270075889
$ws.Cells.Item(19, 13).Value = 1.021049487652812
$ws.Cells.Item(19, 14).Value = 1.012997228400802
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.020552094481935
$ws.Cells.Item(20, 4).Value = 1.032333088153249
$ws.Cells.Item(20, 5).Value = 1.02158104561308
$ws.Cells.Item(20, 6).Value = 1.017027814319169
$ws.Cells.Item(20, 9).Value = 1.032670096706113
$ws.Cells.Item(20, 10).Value = 1.026854053237636
$ws.Cells.Item(20, 11).Value = 1.03571948911965
$ws.Cells.Item(20, 12).Value = 1.025005743028672
$ws.Cells.Item(20, 13).Value = 1.020468998877539
$ws.Cells.Item(20, 14).Value = 1.012895604689852
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.019172748356136
$ws.Cells.Item(21, 4).Value = 1.031656821445924
$ws.Cells.Item(21, 5).Value = 1.020404075245289
$ws.Cells.Item(21, 6).Value = 1.014918523524549
$ws.Cells.Item(21, 9).Value = 1.03238405336151
$ws.Cells.Item(21, 10).Value = 1.025882056438862
$ws.Cells.Item(21, 11).Value = 1.03525319850152
$ws.Cells.Item(21, 12).Value = 1.02404323886449
$ws.Cells.Item(21, 13).Value = 1.018578926833141
$ws.Cells.Item(21, 14).Value = 1.012564231663552
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.018304239060444
$ws.Cells.Item(22, 4).Value = 1.031231241186551
$ws.Cells.Item(22, 5).Value = 1.019663358286711
$ws.Cells.Item(22, 6).Value = 1.013590456162191
$ws.Cells.Item(22, 9).Value = 1.032202538866373
$ws.Cells.Item(22, 10).Value = 1.025269215208335
$ws.Cells.Item(22, 11).Value = 1.034958884873659
$ws.Cells.Item(22, 12).Value = 1.023436739694564
$ws.Cells.Item(22, 13).Value = 1.017388314325815
$ws.Cells.Item(22, 14).Value = 1.012355123382629
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.018764773485177
$ws.Cells.Item(23, 4).Value = 1.031456885799773
$ws.Cells.Item(23, 5).Value = 1.020056094812995
$ws.Cells.Item(23, 6).Value = 1.014294672177904
$ws.Cells.Item(23, 9).Value = 1.032298921524028
$ws.Cells.Item(23, 10).Value = 1.025594257202475
$ws.Cells.Item(23, 11).Value = 1.035115014153234
$ws.Cells.Item(23, 12).Value = 1.02375838407943
$ws.Cells.Item(23, 13).Value = 1.018019698095204
$ws.Cells.Item(23, 14).Value = 1.012466048171605
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.020575198452996
$ws.Cells.Item(24, 4).Value = 1.032344419339101
$ws.Cells.Item(24, 5).Value = 1.021600765936755
$ws.Cells.Item(24, 6).Value = 1.017063146201051
$ws.Cells.Item(24, 9).Value = 1.032674864059216
$ws.Cells.Item(24, 10).Value = 1.026870320231978
$ws.Cells.Item(24, 11).Value = 1.035727287286062
$ws.Cells.Item(24, 12).Value = 1.025021857180175
$ws.Cells.Item(24, 13).Value = 1.020500649024644
$ws.Cells.Item(24, 14).Value = 1.012901147394273
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022670791991497
$ws.Cells.Item(25, 4).Value = 1.033372679120966
$ws.Cells.Item(25, 5).Value = 1.023390308067716
$ws.Cells.Item(25, 6).Value = 1.020268139968773
$ws.Cells.Item(25, 9).Value = 1.033103893678725
$ws.Cells.Item(25, 10).Value = 1.028343815390059
$ws.Cells.Item(25, 11).Value = 1.036432856374802
$ws.Cells.Item(25, 12).Value = 1.026482355108678
$ws.Cells.Item(25, 13).Value = 1.023370295087104
$ws.Cells.Item(25, 14).Value = 1.013402783130343
